$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 380, pushing the existing rows 380-387 down to 382-389.
$ws.Range("A380:T381").EntireRow.Insert()

# --- New row 380: Candy White / Especial ---
$ws.Cells.Item(380, 1).Value = 7
$ws.Cells.Item(380, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(380, 3).Value = "Ñuble"
$ws.Cells.Item(380, 4).Value = 44939
$ws.Cells.Item(380, 5).Value = 16
$ws.Cells.Item(380, 6).Value = "Fruta"
$ws.Cells.Item(380, 7).Value = 100103
$ws.Cells.Item(380, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(380, 9).Value = 100103006
$ws.Cells.Item(380, 10).Value = "Nectarín"
$ws.Cells.Item(380, 11).Value = "Candy White"
$ws.Cells.Item(380, 12).Value = "Especial"
$ws.Cells.Item(380, 13).Value = 80
$ws.Cells.Item(380, 14).Value = 17000
$ws.Cells.Item(380, 15).Value = 17000
$ws.Cells.Item(380, 16).Value = 17000
$ws.Cells.Item(380, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(380, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(380, 19).Value = 1133
$ws.Cells.Item(380, 20).Value = 15

# --- New row 381: Candy White / Primera ---
$ws.Cells.Item(381, 1).Value = 7
$ws.Cells.Item(381, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(381, 3).Value = "Ñuble"
$ws.Cells.Item(381, 4).Value = 44939
$ws.Cells.Item(381, 5).Value = 16
$ws.Cells.Item(381, 6).Value = "Fruta"
$ws.Cells.Item(381, 7).Value = 100103
$ws.Cells.Item(381, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(381, 9).Value = 100103006
$ws.Cells.Item(381, 10).Value = "Nectarín"
$ws.Cells.Item(381, 11).Value = "Candy White"
$ws.Cells.Item(381, 12).Value = "Primera"
$ws.Cells.Item(381, 13).Value = 100
$ws.Cells.Item(381, 14).Value = 15000
$ws.Cells.Item(381, 15).Value = 15000
$ws.Cells.Item(381, 16).Value = 15000
$ws.Cells.Item(381, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(381, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(381, 19).Value = 1000
$ws.Cells.Item(381, 20).Value = 15

# --- Fix the date on the two "Nectar Crest" rows that shifted down (now rows 385-386) ---
$ws.Cells.Item(385, 4).Value = 44249
$ws.Cells.Item(386, 4).Value = 44249

$ws.Range("A1").Select()
